$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ICS2523-3"
$ws.Cells.Item(3, 1).Value = "IEE2463-1"
$ws.Cells.Item(4, 1).Value = "ICE2533-1"
$ws.Cells.Item(5, 1).Value = "ICE3753-1"
$ws.Cells.Item(6, 1).Value = "IIC2764-1"
$ws.Cells.Item(7, 1).Value = "ICE3653-1"
$ws.Cells.Item(8, 1).Value = "IMM2053-1"
$ws.Cells.Item(9, 1).Value = "IIC2333-1"
$ws.Cells.Item(10, 1).Value = "ICC2514-1"
$ws.Cells.Item(11, 1).Value = "IIQ2133-1"
$ws.Cells.Item(12, 1).Value = "ICM2022-1"
$ws.Cells.Item(13, 1).Value = "IIC3143-1"
$ws.Cells.Item(14, 1).Value = "ICC3214-1"
$ws.Cells.Item(15, 1).Value = "ICS3413-1"
$ws.Cells.Item(16, 1).Value = "ICE3124-1"
$ws.Cells.Item(17, 1).Value = "ICE3613-1"
$ws.Cells.Item(18, 1).Value = "IEE2213-1"
$ws.Cells.Item(19, 1).Value = "IMT3150-1"
$ws.Cells.Item(20, 1).Value = "IIQ3343-1"
$ws.Cells.Item(21, 1).Value = "ICE2604-1"
$ws.Cells.Item(22, 1).Value = "IIQ2673-1"
$ws.Cells.Item(23, 1).Value = "IIC2733-1"
$ws.Cells.Item(24, 1).Value = "ICC2105-1"
$ws.Cells.Item(25, 1).Value = "IEE3234-1"
$ws.Cells.Item(26, 1).Value = "ICC3543-1"
$ws.Cells.Item(27, 1).Value = "ICH3364-1"
$ws.Cells.Item(28, 1).Value = "IEE2123-1"
$ws.Cells.Item(29, 1).Value = "ICH2214-1"
$ws.Cells.Item(30, 1).Value = "IIC2733-2"
$ws.Cells.Item(31, 1).Value = "ICS2563-1"
$ws.Cells.Item(32, 1).Value = "ICM2803-1"
$ws.Cells.Item(33, 1).Value = "ICH1005-1"
$ws.Cells.Item(34, 1).Value = "IMM2033-1"
$ws.Cells.Item(35, 1).Value = "ICH3350-1"
$ws.Cells.Item(36, 1).Value = "ICM1001-1"
$ws.Cells.Item(37, 1).Value = "IIC1001-1"
$ws.Cells.Item(38, 1).Value = "ICH3222-1"
$ws.Cells.Item(39, 1).Value = "ICC3124-1"
$ws.Cells.Item(40, 1).Value = "ICM2223-1"
$ws.Cells.Item(41, 1).Value = "ICS2523-4"
$ws.Cells.Item(42, 1).Value = "IIQ2003-1"
$ws.Cells.Item(43, 1).Value = "ICE3233-1"
$ws.Cells.Item(44, 1).Value = "IMM1003-1"
$ws.Cells.Item(45, 1).Value = "IMM2013-1"
$ws.Cells.Item(46, 1).Value = "ICS2123-3"
$ws.Cells.Item(47, 1).Value = "IIC2613-1"
$ws.Cells.Item(48, 1).Value = "ICE2006-1"
$ws.Cells.Item(49, 1).Value = "IMM2213-1"
$ws.Cells.Item(50, 1).Value = "ICE2114-1"
$ws.Cells.Item(51, 1).Value = "ICM2313-1"
$ws.Cells.Item(52, 1).Value = "ICT3523-1"
$ws.Cells.Item(53, 1).Value = "IIC3113-1"
$ws.Cells.Item(54, 1).Value = "IEE2343-1"
$ws.Cells.Item(55, 1).Value = "IIC3113-2"
$ws.Cells.Item(56, 1).Value = "ICT2303-1"
$ws.Cells.Item(57, 1).Value = "IEE2513-1"
$ws.Cells.Item(58, 1).Value = "ICC3253-1"
$ws.Cells.Item(59, 1).Value = "ICM2203-1"
$ws.Cells.Item(60, 1).Value = "ICH2574-1"
$ws.Cells.Item(61, 1).Value = "ICE2020-1"
$ws.Cells.Item(62, 1).Value = "IMM2003-1"
$ws.Cells.Item(63, 1).Value = "IMM3313-1"
$ws.Cells.Item(64, 1).Value = "ICH3532-1"
$ws.Cells.Item(65, 1).Value = "ICE2028-1"
$ws.Cells.Item(66, 1).Value = "ICH3374-1"
$ws.Cells.Item(67, 1).Value = "ICH2304-1"
$ws.Cells.Item(68, 1).Value = "ICH2204-1"
$ws.Cells.Item(69, 1).Value = "ICC3434-1"
$ws.Cells.Item(70, 1).Value = "IEE2113-1"
$ws.Cells.Item(71, 1).Value = "ICM2333-1"
$ws.Cells.Item(72, 1).Value = "ICM3235-1"
$ws.Cells.Item(73, 1).Value = "ICS3313-1"
$ws.Cells.Item(74, 1).Value = "ICM2213-1"
$ws.Cells.Item(75, 1).Value = "IMM2043-1"
$ws.Cells.Item(76, 1).Value = "ICS3811-1"
$ws.Cells.Item(77, 1).Value = "ICC2204-1"
$ws.Cells.Item(78, 1).Value = "ICE2633-1"
$ws.Cells.Item(79, 1).Value = "ICS3723-1"
$ws.Cells.Item(80, 1).Value = "IIQ2303-1"
$ws.Cells.Item(81, 1).Value = "ICE2703-1"
$ws.Cells.Item(82, 1).Value = "ICT3283-1"
$ws.Cells.Item(83, 1).Value = "ICS2123-1"
$ws.Cells.Item(84, 1).Value = "ICE3443-1"
$ws.Cells.Item(85, 1).Value = "ICS3151-1"
$ws.Cells.Item(86, 1).Value = "IMM3323-1"
$ws.Cells.Item(87, 1).Value = "IIC3743-1"
$ws.Cells.Item(88, 1).Value = "IIQ2043-1"
$ws.Cells.Item(89, 1).Value = "IIC2133-1"
$ws.Cells.Item(90, 1).Value = "IIC2133-2"
$ws.Cells.Item(91, 1).Value = "ICE3663-1"
$ws.Cells.Item(92, 1).Value = "ICH2114-1"
$ws.Cells.Item(93, 1).Value = "ICM2403-1"
$ws.Cells.Item(94, 1).Value = "ICE3513-1"
$ws.Cells.Item(95, 1).Value = "IIC3242-1"
$ws.Cells.Item(96, 1).Value = "ICC1001-1"
$ws.Cells.Item(97, 1).Value = "IIC3757-1"
$ws.Cells.Item(98, 1).Value = "ICE2623-1"
$ws.Cells.Item(99, 1).Value = "IEE3732-1"
$ws.Cells.Item(100, 1).Value = "IIC2213-1"
$ws.Cells.Item(101, 1).Value = "ICT3623-1"
$ws.Cells.Item(102, 1).Value = "IMM3800-1"
$ws.Cells.Item(103, 1).Value = "ICH2304-2"
$ws.Cells.Item(104, 1).Value = "ICS2123-2"
$ws.Cells.Item(105, 1).Value = "IEE2713-1"
$ws.Cells.Item(106, 1).Value = "IIC2713-1"
$ws.Cells.Item(107, 1).Value = "ICH3600-1"
$ws.Cells.Item(108, 1).Value = "ICS2023-1"
$ws.Cells.Item(109, 1).Value = "ICM3243-1"
$ws.Cells.Item(110, 1).Value = "ICT3464-1"
$ws.Cells.Item(111, 1).Value = "ICC2304-1"
$ws.Cells.Item(112, 1).Value = "IIC2713-3"
$ws.Cells.Item(113, 1).Value = "IIC2713-2"
$ws.Cells.Item(114, 1).Value = "ICM3762-1"
$ws.Cells.Item(115, 1).Value = "ICS2121-1"
$ws.Cells.Item(116, 1).Value = "ICS3582-1"
$ws.Cells.Item(117, 1).Value = "IMT3800-1"
$ws.Cells.Item(118, 1).Value = "IEE2413-1"
$ws.Cells.Item(119, 1).Value = "IIQ2663-1"
$ws.Cells.Item(120, 1).Value = "ICH2124-1"
$ws.Cells.Item(121, 1).Value = "ICT2233-1"
$ws.Cells.Item(122, 1).Value = "ICS2563-2"
$ws.Cells.Item(123, 1).Value = "ICE3413-1"
$ws.Cells.Item(124, 1).Value = "ICM2003-1"
$ws.Cells.Item(125, 1).Value = "IEE2103-1"
$ws.Cells.Item(126, 1).Value = "IIC2143-2"
$ws.Cells.Item(127, 1).Value = "IIC2143-1"
$ws.Cells.Item(128, 1).Value = "ICM2413-1"
$ws.Cells.Item(129, 1).Value = "ICE2313-1"
$ws.Cells.Item(130, 1).Value = "IEE2613-1"
$ws.Cells.Item(131, 1).Value = "ICT3435-1"
$ws.Cells.Item(132, 1).Value = "IMT3130-1"
$ws.Cells.Item(133, 1).Value = "ICM2028-1"
$ws.Cells.Item(134, 1).Value = "IEE3373-1"
$ws.Cells.Item(135, 1).Value = "IMT2111-1"
$ws.Cells.Item(136, 1).Value = "ICM3251-1"
$ws.Cells.Item(137, 1).Value = "IIC3800-1"
$ws.Cells.Item(138, 1).Value = "ICS3762-1"
$ws.Cells.Item(139, 1).Value = "IIQ3643-1"
$ws.Cells.Item(140, 1).Value = "ICC3264-1"
$ws.Cells.Item(141, 1).Value = "IIC3103-1"
